{"js": "// Update the answer table from the \"before\" problem set to the new\n// (\"c986bee\") problem set. The table has 20 rows x 5 columns; only every\n// 4th row (0, 4, 8, 12, 16) actually holds text - the rows in between are\n// blank rows left for students to write their work. We walk the table in\n// natural reading order (row major, left to right) and, for every cell\n// whose current text matches the next expected \"old\" value, replace it\n// with the corresponding \"new\" value. This keeps the run/paragraph\n// formatting (font, size, justification) untouched because we only set\n// the text of the existing cell body rather than rebuilding the cell.\n\n// Ordered list of (old text -> new text) pairs, in the exact order the\n// cells appear in the document (top-to-bottom, left-to-right).\nconst replacements = [\n  { oldText: \"89\u00f72=44, 1\", newText: \"45\u00f77=6, 3\" },\n  { oldText: \"14\u00f78=1, 6\", newText: \"79\u00f74=19, 3\" },\n  { oldText: \"17\u00f77=2, 3\", newText: \"44\u00f77=6, 2\" },\n  { oldText: \"91\u00f78=11, 3\", newText: \"51\u00f74=12, 3\" },\n  { oldText: \"32\u00f78=4, 0\", newText: \"50\u00f72=25, 0\" },\n  { oldText: \"60\u00f75=12, 0\", newText: \"89\u00f73=29, 2\" },\n  { oldText: \"46\u00f75=9, 1\", newText: \"19\u00f74=4, 3\" },\n  { oldText: \"67\u00f74=16, 3\", newText: \"72\u00f76=12, 0\" },\n  { oldText: \"59\u00f75=11, 4\", newText: \"71\u00f73=23, 2\" },\n  { oldText: \"14\u00f73=4, 2\", newText: \"82\u00f76=13, 4\" },\n  { oldText: \"63\u00f72=31, 1\", newText: \"49\u00f77=7, 0\" },\n  { oldText: \"62\u00f74=15, 2\", newText: \"61\u00f72=30, 1\" },\n  { oldText: \"61\u00f74=15, 1\", newText: \"80\u00f78=10, 0\" },\n  { oldText: \"44\u00f74=11, 0\", newText: \"91\u00f75=18, 1\" },\n  { oldText: \"39\u00f72=19, 1\", newText: \"93\u00f76=15, 3\" },\n  { oldText: \"97\u00f75=19, 2\", newText: \"80\u00f75=16, 0\" },\n  { oldText: \"39\u00f74=9, 3\", newText: \"50\u00f79=5, 5\" },\n  { oldText: \"92\u00f73=30, 2\", newText: \"45\u00f78=5, 5\" },\n  { oldText: \"80\u00f75=16, 0\", newText: \"36\u00f77=5, 1\" },\n  { oldText: \"41\u00f72=20, 1\", newText: \"55\u00f75=11, 0\" },\n  { oldText: \"70\u00f76=11, 4\", newText: \"23\u00f72=11, 1\" },\n  { oldText: \"89\u00f72=44, 1\", newText: \"51\u00f74=12, 3\" },\n  { oldText: \"34\u00f77=4, 6\", newText: \"13\u00f73=4, 1\" },\n  { oldText: \"90\u00f74=22, 2\", newText: \"20\u00f73=6, 2\" },\n  { oldText: \"93\u00f73=31, 0\", newText: \"94\u00f79=10, 4\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in the document body.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Collect every cell in the table (row-major order) so we can walk them\n// in natural reading order and line them up against `replacements`.\nconst rowCount = table.rowCount;\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (let r = 0; r < rowCount; r++) {\n  rows.items[r].cells.load(\"items\");\n}\nawait context.sync();\n\nconst cells = [];\nfor (let r = 0; r < rowCount; r++) {\n  for (const cell of rows.items[r].cells.items) {\n    cells.push(cell);\n    cell.load(\"value\");\n  }\n}\nawait context.sync();\n\n// Setting `cell.value` rewrites only the text of the cell's existing\n// paragraph/run, leaving paragraph and run formatting (font, size,\n// justification, etc.) untouched - unlike `body.insertText(..., replace)`,\n// which would blow away the existing run/paragraph properties.\nlet replacementIndex = 0;\nfor (const cell of cells) {\n  if (replacementIndex >= replacements.length) {\n    break;\n  }\n  const expected = replacements[replacementIndex];\n  if (cell.value === expected.oldText) {\n    cell.value = expected.newText;\n    replacementIndex++;\n  }\n}\n\nawait context.sync();\n\nif (replacementIndex !== replacements.length) {\n  throw new Error(\n    \"Only matched \" +\n      replacementIndex +\n      \" of \" +\n      replacements.length +\n      \" expected cell values; document structure may differ from what was expected.\"\n  );\n}\n", "ps1": "# Update the answer table from the \"before\" problem set to the new\n# (\"c986bee\") problem set. The table has 20 rows x 5 columns; only every\n# 4th row (1, 5, 9, 13, 17 in 1-based COM indexing) actually holds text -\n# the rows in between are blank rows left for students to write their\n# work. We walk the table in natural reading order (row major, left to\n# right) and, for every cell whose current text matches the next expected\n# \"old\" value, replace it with the corresponding \"new\" value via\n# `Cell.Range.Text`, which rewrites only the text inside the existing\n# run/paragraph and therefore preserves formatting (font, size,\n# justification, etc.).\n\n# Ordered list of (old text -> new text) pairs, in the exact order the\n# cells appear in the document (top-to-bottom, left-to-right).\n$replacements = @(\n  @{ Old = \"89\u00f72=44, 1\"; New = \"45\u00f77=6, 3\" },\n  @{ Old = \"14\u00f78=1, 6\"; New = \"79\u00f74=19, 3\" },\n  @{ Old = \"17\u00f77=2, 3\"; New = \"44\u00f77=6, 2\" },\n  @{ Old = \"91\u00f78=11, 3\"; New = \"51\u00f74=12, 3\" },\n  @{ Old = \"32\u00f78=4, 0\"; New = \"50\u00f72=25, 0\" },\n  @{ Old = \"60\u00f75=12, 0\"; New = \"89\u00f73=29, 2\" },\n  @{ Old = \"46\u00f75=9, 1\"; New = \"19\u00f74=4, 3\" },\n  @{ Old = \"67\u00f74=16, 3\"; New = \"72\u00f76=12, 0\" },\n  @{ Old = \"59\u00f75=11, 4\"; New = \"71\u00f73=23, 2\" },\n  @{ Old = \"14\u00f73=4, 2\"; New = \"82\u00f76=13, 4\" },\n  @{ Old = \"63\u00f72=31, 1\"; New = \"49\u00f77=7, 0\" },\n  @{ Old = \"62\u00f74=15, 2\"; New = \"61\u00f72=30, 1\" },\n  @{ Old = \"61\u00f74=15, 1\"; New = \"80\u00f78=10, 0\" },\n  @{ Old = \"44\u00f74=11, 0\"; New = \"91\u00f75=18, 1\" },\n  @{ Old = \"39\u00f72=19, 1\"; New = \"93\u00f76=15, 3\" },\n  @{ Old = \"97\u00f75=19, 2\"; New = \"80\u00f75=16, 0\" },\n  @{ Old = \"39\u00f74=9, 3\"; New = \"50\u00f79=5, 5\" },\n  @{ Old = \"92\u00f73=30, 2\"; New = \"45\u00f78=5, 5\" },\n  @{ Old = \"80\u00f75=16, 0\"; New = \"36\u00f77=5, 1\" },\n  @{ Old = \"41\u00f72=20, 1\"; New = \"55\u00f75=11, 0\" },\n  @{ Old = \"70\u00f76=11, 4\"; New = \"23\u00f72=11, 1\" },\n  @{ Old = \"89\u00f72=44, 1\"; New = \"51\u00f74=12, 3\" },\n  @{ Old = \"34\u00f77=4, 6\"; New = \"13\u00f73=4, 1\" },\n  @{ Old = \"90\u00f74=22, 2\"; New = \"20\u00f73=6, 2\" },\n  @{ Old = \"93\u00f73=31, 0\"; New = \"94\u00f79=10, 4\" }\n)\n\n$d = $word.ActiveDocument\n\nif ($d.Tables.Count -lt 1) {\n  throw \"No table found in the document.\"\n}\n\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$replacementIndex = 0\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    if ($replacementIndex -ge $replacements.Count) {\n      break\n    }\n    $cell = $t.Cell($r, $c)\n    $currentText = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    $expected = $replacements[$replacementIndex]\n    if ($currentText -eq $expected.Old) {\n      $cell.Range.Text = $expected.New\n      $replacementIndex++\n    }\n  }\n}\n\nif ($replacementIndex -ne $replacements.Count) {\n  throw (\"Only matched \" + $replacementIndex + \" of \" + $replacements.Count + \" expected cell values; document structure may differ from what was expected.\")\n}\n"}
